{"js": "const replacements = [\n  [\"2025-05-23 Friday\", \"2025-05-24 Saturday\"],\n  [\"18\u00d792=1656\", \"35\u00d762=2170\"],\n  [\"70\u00d766=4620\", \"74\u00d765=4810\"],\n  [\"70\u00d767=4690\", \"87\u00d738=3306\"],\n  [\"15\u00d713=195\", \"18\u00d713=234\"],\n  [\"34\u00d765=2210\", \"86\u00d774=6364\"],\n  [\"50\u00d775=3750\", \"69\u00d742=2898\"],\n  [\"62\u00d757=3534\", \"51\u00d771=3621\"],\n  [\"45\u00d771=3195\", \"49\u00d764=3136\"],\n  [\"55\u00d764=3520\", \"79\u00d765=5135\"],\n  [\"52\u00d776=3952\", \"31\u00d773=2263\"],\n  [\"70\u00d780=5600\", \"73\u00d724=1752\"],\n  [\"82\u00d759=4838\", \"71\u00d791=6461\"],\n  [\"90\u00d794=8460\", \"75\u00d720=1500\"],\n  [\"67\u00d737=2479\", \"48\u00d778=3744\"],\n  [\"48\u00d740=1920\", \"52\u00d722=1144\"],\n  [\"24\u00d744=1056\", \"12\u00d735=420\"],\n  [\"34\u00d767=2278\", \"82\u00d796=7872\"],\n  [\"92\u00d714=1288\", \"45\u00d724=1080\"],\n  [\"75\u00d713=975\", \"20\u00d719=380\"],\n  [\"44\u00d719=836\", \"14\u00d732=448\"],\n  [\"13\u00d764=832\", \"61\u00d726=1586\"],\n  [\"59\u00d746=2714\", \"99\u00d740=3960\"],\n  [\"73\u00d748=3504\", \"40\u00d734=1360\"],\n  [\"75\u00d739=2925\", \"83\u00d757=4731\"],\n  [\"90\u00d720=1800\", \"97\u00d717=1649\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-05-23 Friday\", \"2025-05-24 Saturday\"),\n    @(\"18\u00d792=1656\", \"35\u00d762=2170\"),\n    @(\"70\u00d766=4620\", \"74\u00d765=4810\"),\n    @(\"70\u00d767=4690\", \"87\u00d738=3306\"),\n    @(\"15\u00d713=195\", \"18\u00d713=234\"),\n    @(\"34\u00d765=2210\", \"86\u00d774=6364\"),\n    @(\"50\u00d775=3750\", \"69\u00d742=2898\"),\n    @(\"62\u00d757=3534\", \"51\u00d771=3621\"),\n    @(\"45\u00d771=3195\", \"49\u00d764=3136\"),\n    @(\"55\u00d764=3520\", \"79\u00d765=5135\"),\n    @(\"52\u00d776=3952\", \"31\u00d773=2263\"),\n    @(\"70\u00d780=5600\", \"73\u00d724=1752\"),\n    @(\"82\u00d759=4838\", \"71\u00d791=6461\"),\n    @(\"90\u00d794=8460\", \"75\u00d720=1500\"),\n    @(\"67\u00d737=2479\", \"48\u00d778=3744\"),\n    @(\"48\u00d740=1920\", \"52\u00d722=1144\"),\n    @(\"24\u00d744=1056\", \"12\u00d735=420\"),\n    @(\"34\u00d767=2278\", \"82\u00d796=7872\"),\n    @(\"92\u00d714=1288\", \"45\u00d724=1080\"),\n    @(\"75\u00d713=975\", \"20\u00d719=380\"),\n    @(\"44\u00d719=836\", \"14\u00d732=448\"),\n    @(\"13\u00d764=832\", \"61\u00d726=1586\"),\n    @(\"59\u00d746=2714\", \"99\u00d740=3960\"),\n    @(\"73\u00d748=3504\", \"40\u00d734=1360\"),\n    @(\"75\u00d739=2925\", \"83\u00d757=4731\"),\n    @(\"90\u00d720=1800\", \"97\u00d717=1649\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 0\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute([ref]$null, $false, $false, $false, $false, $false, $true, 0, $false, $newText, 2)\n}\n"}
